$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 34 (pushes existing rows 34-90 down to 35-91).
$ws.Rows.Item(34).Insert()

# Populate the newly-inserted row 34 with the new weekly record.
$ws.Cells.Item(34, 1).Value = 11
$ws.Cells.Item(34, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(34, 3).Value = "Bíobío"
$ws.Cells.Item(34, 4).Value = 45100
$ws.Cells.Item(34, 5).Value = 8
$ws.Cells.Item(34, 6).Value = 100112013
$ws.Cells.Item(34, 7).Value = "Alcachofa"
$ws.Cells.Item(34, 8).Value = "Argentina(o)"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 110
$ws.Cells.Item(34, 11).Value = 14000
$ws.Cells.Item(34, 12).Value = 15000
$ws.Cells.Item(34, 13).Value = 14455
$ws.Cells.Item(34, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 289
$ws.Cells.Item(34, 17).Value = 50
$ws.Cells.Item(34, 18).Value = "Hortaliza"
